# Verify_40V_On_Addition_Deletion_Of_Rbus.xlsx
# "Updated test data as per new implemenation"

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Add Panels")

# 1. Loading-details label text changed ("40V (A)" -> "40V Rail(A)") for the
#    whole data column (F8:F12 all shared the same string).
$ws1.Range("F8:F12").Value = "40V Rail(A)"

# 2. Add the new, still-empty "Test data" worksheet right after "Add Panels".
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "Test data"
$newSheet.Range("A1:F3").Select() | Out-Null

# 3. Re-activate "Add Panels" (stays the selected/visible tab) and leave it
#    scrolled/selected on the new range.
$ws1.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 3
$ws1.Range("A10:F12").Select() | Out-Null
